$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.931.56'
$ws.Range("E2").Value = '  +4.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.738.97'
$ws.Range("E3").Value = '  +3.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.81'
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.17'
$ws.Range("E6").Value = '  +10.19%  '

$ws.Range("E7").Value = '  +3.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.758.72'
$ws.Range("E9").Value = '  +3.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.82'
$ws.Range("E10").Value = '  +3.39%  '

$ws.Range("E11").Value = '  +1.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.395'
$ws.Range("E12").Value = '  +3.48%  '

$ws.Range("E13").Value = '  +0.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.225.95'
$ws.Range("E14").Value = '  +2.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.71'
$ws.Range("E15").Value = '  +5.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.853.34'
$ws.Range("E16").Value = '  +4.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000156'
$ws.Range("E17").Value = '  +6.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.752.94'
$ws.Range("E18").Value = '  +3.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.21'
$ws.Range("E19").Value = '  +4.40%  '

$ws.Range("E20").Value = '  +3.72%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '363.62'
$ws.Range("E21").Value = '  +2.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.01'
$ws.Range("E22").Value = '  +1.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.542'
$ws.Range("E23").Value = '  +2.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.993'
$ws.Range("E24").Value = '  -0.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.89'
$ws.Range("E25").Value = '  +3.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.172'
$ws.Range("E26").Value = '  +5.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.60'
$ws.Range("E27").Value = '  +0.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0921'
$ws.Range("E29").Value = '  +11.78%  '

$ws.Range("E30").Value = '  +1.13%  '

$ws.Range("E31").Value = '  +4.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.26'
$ws.Range("E32").Value = '  +12.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '173.41'
$ws.Range("E33").Value = '  +2.72%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.63'
$ws.Range("E35").Value = '  +2.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.99'
$ws.Range("E36").Value = '  +6.14%  '

$ws.Range("E37").Value = '  +5.90%  '

$ws.Range("E38").Value = '  +5.54%  '

$ws.Range("E39").Value = '  +3.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.25'
$ws.Range("E40").Value = '  +1.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '338.07'
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.20'
$ws.Range("E42").Value = '  +16.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.56'
$ws.Range("E43").Value = '  +2.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.46'
$ws.Range("E44").Value = '  +6.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.86'
$ws.Range("E45").Value = '  +5.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0603'
$ws.Range("E46").Value = '  +3.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.645'
$ws.Range("E47").Value = '  +2.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0259'
$ws.Range("E48").Value = '  +2.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.37'
$ws.Range("E49").Value = '  +1.20%  '

$ws.Range("E50").Value = '  +2.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.17%  '
